# Adding Budget sonar results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: rename label "sonar" -> "sonar 157/51" and fill in the Budget
# (KPB) training/test accuracy results that were previously blank.
$ws.Range("A4").Value = "sonar 157/51"
$ws.Range("E4").Value = 0.5414
$ws.Range("F4").Value = 0.5098

# Row 5: rename label "sonar" -> "sonar 116/92" and fill in the Budget
# (KPB) training/test accuracy results that were previously blank.
$ws.Range("A5").Value = "sonar 116/92"
$ws.Range("E5").Value = 0.5517
$ws.Range("F5").Value = 0.5109
